# Weekly update: insert the newest week (2022-10-17, serial 44841) of
# "Chirimoya" price data at the top of the data block (rows 164-169),
# pushing the existing history down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 164..184 down by 6 to make room for the new week's 6 rows.
$ws.Range("A164:T169").EntireRow.Insert()

# New week's data (date serial 44841), one row per "Calidad" grade,
# in the same order the sheet already uses for this market/week block.
# (Single-quoted literals below so the '$' in "$/kilo" / "$/bandeja" is
# never treated as PowerShell variable interpolation.)
$newRows = @(
    @{ L = 'Cuarta';                  M = 40;  N = 1200;  O = 1200;  P = 1200;  Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia del Elquí'; S = 1200; T = 1 }
    @{ L = 'Especial';                M = 75;  N = 22400; O = 22400; P = 22400; Q = '$/bandeja 8 kilos';            R = 'Provincia del Elquí'; S = 2800; T = 8 }
    @{ L = 'Extra (doble especial)';  M = 50;  N = 24000; O = 24000; P = 24000; Q = '$/bandeja 8 kilos';            R = 'Provincia del Elquí'; S = 3000; T = 8 }
    @{ L = 'Primera';                 M = 110; N = 19200; O = 19200; P = 19200; Q = '$/bandeja 8 kilos';            R = 'Provincia del Elquí'; S = 2400; T = 8 }
    @{ L = 'Segunda';                 M = 60;  N = 2200;  O = 2200;  P = 2200;  Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia del Elquí'; S = 2200; T = 1 }
    @{ L = 'Tercera';                 M = 55;  N = 1600;  O = 1600;  P = 1600;  Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia del Elquí'; S = 1600; T = 1 }
)

$r = 164
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = 'Vega Central Mapocho de Santiago'
    $ws.Cells.Item($r, 3).Value = 'Metropolitana'
    $ws.Cells.Item($r, 4).Value = 44841
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 'Fruta'
    $ws.Cells.Item($r, 7).Value = 100107
    $ws.Cells.Item($r, 8).Value = 'Otros'
    $ws.Cells.Item($r, 9).Value = 100107002
    $ws.Cells.Item($r, 10).Value = 'Chirimoya'
    $ws.Cells.Item($r, 11).Value = 'Cultivar IV Región'
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r = $r + 1
}
